$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Unprotect the sheet so we can update protected cells
$ws.Unprotect("D382")

# Update the confidential disclaimer date from 2021-05-13 to 2021-05-14
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-38
$ws.Range("D2").Value = 0.03265319185750534
$ws.Range("E2").Value = -0.001457975986277837
$ws.Range("D3").Value = 0.02845412420574872
$ws.Range("E3").Value = 0.02615151003880545
$ws.Range("D4").Value = 0.02777861680007121
$ws.Range("E4").Value = 0.02052264331645914
$ws.Range("D5").Value = 0.06323945507335667
$ws.Range("E5").Value = 0.01943083439033111
$ws.Range("D6").Value = 0.01605255235582458
$ws.Range("E6").Value = 0.003214953271028165
$ws.Range("D7").Value = 0.01525562564447747
$ws.Range("E7").Value = 0.02391629297458908
$ws.Range("D8").Value = 0.03025321026043141
$ws.Range("E8").Value = -0.00273733486729888
$ws.Range("D9").Value = 0.03416943293983744
$ws.Range("E9").Value = 0.02423603793466822
$ws.Range("D10").Value = 0.02920424347471248
$ws.Range("E10").Value = 0.01378786010767286
$ws.Range("D11").Value = 0.02885578807030318
$ws.Range("E11").Value = 0.01664402173913038
$ws.Range("D12").Value = 0.01082572128968506
$ws.Range("E12").Value = 0.04693274205469322
$ws.Range("D13").Value = 0.01418145485132136
$ws.Range("E13").Value = 0.02172195892575046
$ws.Range("D14").Value = 0.01402182947088587
$ws.Range("E14").Value = 0.04563610944677432
$ws.Range("D15").Value = 0.009179859597676032
$ws.Range("E15").Value = 0.01821668264621268
$ws.Range("D16").Value = 0.008149895833437518
$ws.Range("E16").Value = 0.03379721669980107
$ws.Range("D17").Value = 0.02966511675356382
$ws.Range("E17").Value = 0.02737656943264422
$ws.Range("D18").Value = 0.02568248351487887
$ws.Range("E18").Value = 0.008481836874571602
$ws.Range("D19").Value = 0.0333172974502946
$ws.Range("E19").Value = 0.01909221902017277
$ws.Range("D20").Value = 0.03053085440585053
$ws.Range("E20").Value = 0.03498656882657403
$ws.Range("D21").Value = 0.04524659420847914
$ws.Range("E21").Value = 0.02395699324040557
$ws.Range("D22").Value = 0.03588670598625181
$ws.Range("E22").Value = 0.02775842367826975
$ws.Range("D23").Value = 0.03254717500082514
$ws.Range("E23").Value = -0.005500583860856723
$ws.Range("D24").Value = 0.03168323763478394
$ws.Range("E24").Value = 0.002841070515370081
$ws.Range("D25").Value = 0.01412504588229529
$ws.Range("E25").Value = 0.07405047157787403
$ws.Range("D26").Value = 0.01469273614504706
$ws.Range("E26").Value = 0.03471654958340142
$ws.Range("D27").Value = 0.0317224438685751
$ws.Range("E27").Value = 0.003537490462648396
$ws.Range("D28").Value = 0.0313265809263673
$ws.Range("E28").Value = -0.003882304863097596
$ws.Range("D29").Value = 0.02916823774980222
$ws.Range("E29").Value = 0.021067357939349
$ws.Range("D30").Value = 0.02939887442103295
$ws.Range("E30").Value = 0.01736396976274235
$ws.Range("D31").Value = 0.0328018154886627
$ws.Range("E31").Value = 0.04227877279961967
$ws.Range("D32").Value = 0.03220652083681306
$ws.Range("E32").Value = 0.001503040240486309
$ws.Range("D33").Value = 0.02827949643993395
$ws.Range("E33").Value = 0.0341149425287357
$ws.Range("D34").Value = 0.03268999770963583
$ws.Range("E34").Value = 0.003035049931466638
$ws.Range("D35").Value = 0.03106573945257296
$ws.Range("E35").Value = 0
$ws.Range("D36").Value = 0.03230513651670616
$ws.Range("E36").Value = 0.01554179566563452
$ws.Range("D37").Value = 0.0333829078823533
$ws.Range("E37").Value = 0.01222376683763948
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0.01724514197757432

# Re-protect the sheet with the original password
$ws.Protect("D382")
